$d = $word.ActiveDocument

# Locate the "END OF SECTION" paragraph (it currently sits right after the
# section title, near the top of the document) and remove it entirely,
# including its paragraph mark.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $paraText = $para.Range.Text.TrimEnd([char]13)
    if ($paraText -eq "END OF SECTION") {
        $para.Range.Delete()
        break
    }
}

# Re-append "END OF SECTION" as the final paragraph of the document body,
# right before the section break, matching where it belongs at the end.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.InsertParagraphAfter()

$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newIndex)
$newPara.Range.Text = "END OF SECTION"
